$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    12  = -6.764
    32  = -7.276000000000001
    36  = -7.797000000000001
    38  = -7.771999999999998
    46  = -8.156000000000001
    54  = -7.877000000000001
    55  = -8.028
    67  = -7.557
    69  = -7.737
    72  = -7.398000000000001
    91  = -7.636
    99  = -8.253
    104 = -7.806999999999999
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 4).Value = $updates[$row]
}
